# Adds row 9 (ATGE / Adtalem Global Education Inc.) to Sheet1
# and extends the used range dimension from A1:II8 to A1:II9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = 9

$ws.Range("A$r").Value = "ATGE"
$ws.Range("B$r").Value = "Education & Training Services"
$ws.Range("C$r").Value = "US"
$ws.Range("D$r").Value = "Adtalem Global Education Inc."
$ws.Range("E$r").Value = 28.48
$ws.Range("F$r").Value = 7
$ws.Range("G$r").Value = 26.66
$ws.Range("H$r").Value = -35
$ws.Range("I$r").Value = 43.85
$ws.Range("J$r").Value = 50
$ws.Range("K$r").Value = 87
$ws.Range("L$r").Value = 0
$ws.Range("M$r").Value = 1
$ws.Range("N$r").Value = 24.24
$ws.Range("O$r").Value = 94.51000000000001
$ws.Range("P$r").Value = -14.93
$ws.Range("Q$r").Value = 5.74
$ws.Range("R$r").Value = 7.44
$ws.Range("S$r").Value = 78.47
$ws.Range("T$r").Value = 0.89
$ws.Range("U$r").Value = -0.29
$ws.Range("V$r").Value = 0
$ws.Range("W$r").Value = 0.46
$ws.Range("X$r").Value = 730464
$ws.Range("Y$r").Value = 1417298688
$ws.Range("Z$r").Value = "'False"
$ws.Range("Z$r").ClearFormats()
$ws.Range("AA$r").Value = 12585
$ws.Range("AB$r").Value = "'2021-06-30"
$ws.Range("AB$r").ClearFormats()
$ws.Range("AC$r").Value = "USD"
$ws.Range("AD$r").Value = "'2021-08-19"
$ws.Range("AD$r").ClearFormats()
$ws.Range("AE$r").Value = "'2021-08-19"
$ws.Range("AE$r").ClearFormats()
$ws.Range("AF$r").Value = 2021
$ws.Range("AG$r").Value = "FY"
$ws.Range("AH$r").Value = 1112380000
$ws.Range("AI$r").Value = 489233000
$ws.Range("AJ$r").Value = 623147000
$ws.Range("AK$r").Value = 1
$ws.Range("AL$r").Value = 0
$ws.Range("AM$r").Value = 0
$ws.Range("AN$r").Value = 0
$ws.Range("AO$r").Value = 420267000
$ws.Range("AP$r").Value = 0
$ws.Range("AQ$r").Value = 420267000
$ws.Range("AR$r").Value = 909500000
$ws.Range("AS$r").Value = 4094000
$ws.Range("AT$r").Value = 41365000
$ws.Range("AU$r").Value = 99051000
$ws.Range("AV$r").Value = 242573000
$ws.Range("AW$r").Value = 0
$ws.Range("AX$r").Value = 202880000
$ws.Range("AY$r").Value = 0
$ws.Range("AZ$r").Value = 76030000
$ws.Range("BA$r").Value = 126850000
$ws.Range("BB$r").Value = 0
$ws.Range("BC$r").Value = 25248000
$ws.Range("BD$r").Value = 76909000
$ws.Range("BE$r").Value = 0
$ws.Range("BF$r").Value = 2
$ws.Range("BG$r").Value = 1
$ws.Range("BH$r").Value = 51322000
$ws.Range("BI$r").Value = 51645000
$ws.Range("BJ$r").Value = "https://www.sec.gov/Archives/edgar/data/730464/000155837021011957/0001558370-21-011957-index.htm"
$ws.Range("BK$r").Value = "https://www.sec.gov/Archives/edgar/data/730464/000155837021011957/atge-20210630x10k.htm"
$ws.Range("BL$r").Value = 494613000
$ws.Range("BM$r").Value = 0
$ws.Range("BN$r").Value = 494613000
$ws.Range("BO$r").Value = 67996000
$ws.Range("BP$r").Value = 0
$ws.Range("BQ$r").Value = 952344000
$ws.Range("BR$r").Value = 1514953000
$ws.Range("BS$r").Value = 466180000
$ws.Range("BT$r").Value = 686374000
$ws.Range("BU$r").Value = 276249000
$ws.Range("BV$r").Value = 962623000
$ws.Range("BW$r").Value = 0
$ws.Range("BX$r").Value = 22479000
$ws.Range("BY$r").Value = 87601000
$ws.Range("BZ$r").Value = 1538883000
$ws.Range("CA$r").Value = 0
$ws.Range("CB$r").Value = 3053836000
$ws.Range("CC$r").Value = 56071000
$ws.Range("CD$r").Value = 58329000
$ws.Range("CE$r").Value = 0
$ws.Range("CF$r").Value = 100697000
$ws.Range("CG$r").Value = 193710000
$ws.Range("CH$r").Value = 408807000
$ws.Range("CI$r").Value = 1235566000
$ws.Range("CJ$r").Value = 0
$ws.Range("CK$r").Value = 26991000
$ws.Range("CL$r").Value = 81402000
$ws.Range("CM$r").Value = 1343959000
$ws.Range("CN$r").Value = 0
$ws.Range("CO$r").Value = 223184000
$ws.Range("CP$r").Value = 1752766000
$ws.Range("CQ$r").Value = 0
$ws.Range("CR$r").Value = 811000
$ws.Range("CS$r").Value = 2005105000
$ws.Range("CT$r").Value = -7365000
$ws.Range("CU$r").Value = -697481000
$ws.Range("CV$r").Value = 1301070000
$ws.Range("CW$r").Value = 0
$ws.Range("CX$r").Value = 1301070000
$ws.Range("CY$r").Value = 3053836000
$ws.Range("CZ$r").Value = 3053836000
$ws.Range("DA$r").Value = 0
$ws.Range("DB$r").Value = 1293895000
$ws.Range("DC$r").Value = 799282000
$ws.Range("DD$r").Value = 1519000
$ws.Range("DE$r").Value = 13875000
$ws.Range("DF$r").Value = -7478000
$ws.Range("DG$r").Value = 13259000
$ws.Range("DH$r").Value = 8530000
$ws.Range("DI$r").Value = -6638000
$ws.Range("DJ$r").Value = 8757000
$ws.Range("DK$r").Value = 192199000
$ws.Range("DL$r").Value = 0
$ws.Range("DM$r").Value = 0
$ws.Range("DN$r").Value = -10745000
$ws.Range("DO$r").Value = 2721000
$ws.Range("DP$r").Value = -48664000
$ws.Range("DQ$r").Value = -56688000
$ws.Range("DR$r").Value = -3000000
$ws.Range("DS$r").Value = 0
$ws.Range("DT$r").Value = -100000000
$ws.Range("DU$r").Value = 0
$ws.Range("DV$r").Value = 779466000
$ws.Range("DW$r").Value = 676466000
$ws.Range("DX$r").Value = 534000
$ws.Range("DY$r").Value = 812511000
$ws.Range("DZ$r").Value = 1313616000
$ws.Range("EA$r").Value = 501105000
$ws.Range("EB$r").Value = 192199000
$ws.Range("EC$r").Value = -48664000
$ws.Range("ED$r").Value = 143535000
$ws.Range("EE$r").Value = 1289015800
$ws.Range("EF$r").Value = 642125000
$ws.Range("EG$r").Value = 433876000
$ws.Range("EH$r").Value = "Adtalem Global Education Inc."
$ws.Range("EI$r").Value = -0
$ws.Range("EJ$r").Value = -0
$ws.Range("EK$r").Value = 28
$ws.Range("EL$r").Value = 29
$ws.Range("EM$r").Value = 30
$ws.Range("EN$r").Value = 35
$ws.Range("EO$r").Value = 274479
$ws.Range("EP$r").Value = 399448
$ws.Range("EQ$r").Value = "NYSE"
$ws.Range("ER$r").Value = 29
$ws.Range("ES$r").Value = 29
$ws.Range("ET$r").Value = 0
$ws.Range("EU$r").Value = 49764701
$ws.Range("EV$r").Value = 1643525203
$ws.Range("EW$r").Value = 1052001000
$ws.Range("EX$r").Value = 395838000
$ws.Range("EY$r").Value = 692766000
$ws.Range("EZ$r").Value = 91589000
$ws.Range("FA$r").Value = 345983000
$ws.Range("FB$r").Value = 1310421000
$ws.Range("FC$r").Value = 16275000
$ws.Range("FD$r").Value = 107692000
$ws.Range("FE$r").Value = 1239687000
$ws.Range("FF$r").Value = 400411000
$ws.Range("FG$r").Value = 504700000
$ws.Range("FH$r").Value = 99790000
$ws.Range("FI$r").Value = 311631000
$ws.Range("FJ$r").Value = 1391530000
$ws.Range("FK$r").Value = 348327000
$ws.Range("FL$r").Value = 159479000
$ws.Range("FM$r").Value = 709257000
$ws.Range("FN$r").Value = 266654000
$ws.Range("FO$r").Value = 635695000
$ws.Range("FP$r").Value = 1265181000
$ws.Range("FQ$r").Value = 1828317000
$ws.Range("FR$r").Value = 1468222000
$ws.Range("FS$r").Value = 40864000
$ws.Range("FT$r").Value = 280374000
$ws.Range("FU$r").Value = 107868000
$ws.Range("FV$r").Value = 1514953000
$ws.Range("FW$r").Value = 100697000
$ws.Range("FX$r").Value = 408807000
$ws.Range("FY$r").Value = 1301070000
$ws.Range("FZ$r").Value = 1293895000
$ws.Range("GA$r").Value = 799282000
$ws.Range("GB$r").Value = 48036000
$ws.Range("GC$r").Value = 280654000
$ws.Range("GD$r").Value = 108500000
$ws.Range("GE$r").Value = 1497278000
$ws.Range("GF$r").Value = 116670000
$ws.Range("GG$r").Value = 382205000
$ws.Range("GH$r").Value = 1306224000
$ws.Range("GI$r").Value = 792115000
$ws.Range("GJ$r").Value = 80770000
$ws.Range("GK$r").Value = 1192466000
$ws.Range("GL$r").Value = 153136000
$ws.Range("GM$r").Value = -41545000
$ws.Range("GN$r").Value = 31
$ws.Range("GO$r").Value = 1
$ws.Range("GP$r").Value = 399032
$ws.Range("GQ$r").Value = 1526532096
$ws.Range("GR$r").Value = 0
$ws.Range("GS$r").Value = "26.66-43.85"
$ws.Range("GT$r").Value = 0
$ws.Range("GU$r").Value = "USD"
$ws.Range("GV$r").Value = "US00737L1035"
$ws.Range("GW$r").Value = "00737L103"
$ws.Range("GX$r").Value = "NYSE"
$ws.Range("GY$r").Value = "https://www.adtalem.com"
$ws.Range("GZ$r").Value = "Adtalem Global Education Inc. provides workforce solutions worldwide. It operates through two segments, Medical and Healthcare; and Financial Services. The Medical and Healthcare segment offers degree and non-degree programs in the medical and healthcare postsecondary education industry. This segment operates Chamberlain University, American University of the Caribbean School of Medicine, Ross University School of Medicine, and Ross University School of Veterinary Medicine. The Financial Services segment provides test preparation, certifications, conferences, seminars, memberships, and subscriptions to business professionals in the areas of accounting, anti-money laundering, banking, and mortgage lending. It operates Association of Certified Anti-Money Laundering Specialists, Becker Professional Education, OnCourse Learning, and EduPristine. The company was formerly known as DeVry Education Group Inc. and changed its name to Adtalem Global Education Inc. in May 2017. Adtalem Global Education Inc. was incorporated in 1987 and is based in Chicago, Illinois."
$ws.Range("HA$r").Value = "Ms. Lisa Wardell"
$ws.Range("HB$r").Value = "Consumer Defensive"
$ws.Range("HC$r").Value = 4426
$ws.Range("HD$r").Value = "'16305157700"
$ws.Range("HD$r").ClearFormats()
$ws.Range("HE$r").Value = "500 W Monroe St Fl 28"
$ws.Range("HF$r").Value = "Chicago"
$ws.Range("HG$r").Value = "ILLINOIS"
$ws.Range("HH$r").Value = "'60661"
$ws.Range("HH$r").ClearFormats()
$ws.Range("HI$r").Value = 0
$ws.Range("HJ$r").Value = 40
$ws.Range("HK$r").Value = "https://fmpcloud.io/image-stock/ATGE.png"
$ws.Range("HL$r").Value = "'1995-11-14"
$ws.Range("HL$r").ClearFormats()
$ws.Range("HM$r").Value = "'False"
$ws.Range("HM$r").ClearFormats()
$ws.Range("HN$r").Value = "'False"
$ws.Range("HN$r").ClearFormats()
$ws.Range("HO$r").Value = "'True"
$ws.Range("HO$r").ClearFormats()
$ws.Range("HP$r").Value = "'False"
$ws.Range("HP$r").ClearFormats()
$ws.Range("HQ$r").Value = 9108000
$ws.Range("HR$r").Value = -8201000
$ws.Range("HS$r").Value = 1121488000
$ws.Range("HT$r").Value = 1043800000
$ws.Range("HU$r").Value = 165957000
$ws.Range("HV$r").Value = -15973000
$ws.Range("HW$r").Value = 514284000
$ws.Range("HX$r").Value = 264401000
$ws.Range("HY$r").Value = 80086000
$ws.Range("HZ$r").Value = 0
$ws.Range("IA$r").Value = 26956530
$ws.Range("IB$r").Value = -68501530
$ws.Range("IC$r").Value = 25
$ws.Range("ID$r").Value = 84634470
$ws.Range("IE$r").Value = 2
$ws.Range("IF$r").Value = 73562000
$ws.Range("IG$r").Value = 1
$ws.Range("IH$r").Value = 0.05
$ws.Range("II$r").Value = -411018312
